# amazonBot ASIN tracking sheet update:
#  - replace the leftover "test" placeholder in row 8 with the real
#    product label ("5600x"), keeping its existing ASIN/price
#  - append a new tracked ASIN row (5800x) used while testing the new
#    duplicate-order check / Discord webhook-on-checkout logic
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = "5600x"
$ws.Range("B8").Value = "B08166SLDF"
$ws.Range("C8").Value = 300

$ws.Range("A9").Value = "5800x"
$ws.Range("B9").Value = "B0815XFSGK"
$ws.Range("C9").Value = 425

$ws.Range("C11").Select()
